# Insert a new weekly record for "Feria Lagunitas de Puerto Montt - Cebollín"
# at row 174, pushing the existing rows 174:226 down to 175:227.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 174 (shifts 174..226 -> 175..227).
$ws.Rows.Item(174).Insert()

# Populate the newly inserted row with the new weekly data point.
$ws.Range("A174").Value = 4
$ws.Range("B174").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C174").Value = "Los Lagos"
$ws.Range("D174").Value = 44588
$ws.Range("E174").Value = 10
$ws.Range("F174").Value = 100112037
$ws.Range("G174").Value = "Cebollín"
$ws.Range("H174").Value = "Sin especificar"
$ws.Range("I174").Value = "Primera"
$ws.Range("J174").Value = 70
$ws.Range("K174").Value = 6000
$ws.Range("L174").Value = 6500
$ws.Range("M174").Value = 6250
$ws.Range("N174").Value = "$/paquete 36 unidades"
$ws.Range("O174").Value = "Región Metropolitana"
$ws.Range("P174").Value = 174
$ws.Range("Q174").Value = 36
$ws.Range("R174").Value = "Hortaliza"
